$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value2 = 4
$ws.Range("F2").Value2 = 163
$ws.Range("H2").Value2 = "bedrooms"
$ws.Range("L2").Value2 = "stimuli/img_x0u5z.png"
$ws.Range("M2").Value2 = 92
$ws.Range("N2").Value2 = 78.16216216216216
$ws.Range("O2").Value2 = 85.08108108108108
$ws.Range("P2").Value2 = 37
$ws.Range("C3").Value2 = 4
$ws.Range("F3").Value2 = 164
$ws.Range("H3").Value2 = "kitchens"
$ws.Range("L3").Value2 = "stimuli/img_uegbb.png"
$ws.Range("M3").Value2 = 78.80952380952381
$ws.Range("N3").Value2 = 61.52380952380953
$ws.Range("O3").Value2 = 70.16666666666667
$ws.Range("Q3").Value2 = 8
$ws.Range("R3").Value2 = 8
$ws.Range("S3").Value2 = 8
$ws.Range("C4").Value2 = 4
$ws.Range("F4").Value2 = 165
$ws.Range("H4").Value2 = "bedrooms"
$ws.Range("L4").Value2 = "stimuli/img_t2ioc.png"
$ws.Range("M4").Value2 = 88.18918918918919
$ws.Range("N4").Value2 = 74.05405405405405
$ws.Range("O4").Value2 = 81.12162162162161
$ws.Range("P4").Value2 = 37
$ws.Range("C5").Value2 = 4
$ws.Range("F5").Value2 = 166
$ws.Range("H5").Value2 = "bedrooms"
$ws.Range("L5").Value2 = "stimuli/img_le8uf.png"
$ws.Range("M5").Value2 = 12.88888888888889
$ws.Range("N5").Value2 = 9.222222222222221
$ws.Range("O5").Value2 = 11.05555555555556
$ws.Range("P5").Value2 = 36
$ws.Range("C6").Value2 = 4
$ws.Range("F6").Value2 = 167
$ws.Range("H6").Value2 = "bedrooms"
$ws.Range("L6").Value2 = "stimuli/img_rvssl.png"
$ws.Range("M6").Value2 = 74.25
$ws.Range("N6").Value2 = 54.33333333333334
$ws.Range("O6").Value2 = 64.29166666666667
$ws.Range("P6").Value2 = 36
$ws.Range("C7").Value2 = 4
$ws.Range("F7").Value2 = 168
$ws.Range("H7").Value2 = "bedrooms"
$ws.Range("L7").Value2 = "stimuli/img_bj2gr.png"
$ws.Range("M7").Value2 = 65.25
$ws.Range("N7").Value2 = 44.8
$ws.Range("O7").Value2 = 55.025
$ws.Range("P7").Value2 = 40
$ws.Range("C8").Value2 = 4
$ws.Range("F8").Value2 = 169
$ws.Range("H8").Value2 = "bedrooms"
$ws.Range("L8").Value2 = "stimuli/img_qgbyn.png"
$ws.Range("M8").Value2 = 65.08108108108108
$ws.Range("N8").Value2 = 40.10810810810811
$ws.Range("O8").Value2 = 52.5945945945946
$ws.Range("P8").Value2 = 37
$ws.Range("C9").Value2 = 4
$ws.Range("F9").Value2 = 170
$ws.Range("H9").Value2 = "bedrooms"
$ws.Range("L9").Value2 = "stimuli/img_h0hbk.png"
$ws.Range("M9").Value2 = 86.80952380952381
$ws.Range("N9").Value2 = 69.19047619047619
$ws.Range("O9").Value2 = 78
$ws.Range("P9").Value2 = 42
$ws.Range("C10").Value2 = 4
$ws.Range("F10").Value2 = 171
$ws.Range("H10").Value2 = "bedrooms"
$ws.Range("L10").Value2 = "stimuli/img_oou46.png"
$ws.Range("M10").Value2 = 75.70270270270271
$ws.Range("N10").Value2 = 54.86486486486486
$ws.Range("O10").Value2 = 65.28378378378379
$ws.Range("P10").Value2 = 37
$ws.Range("C11").Value2 = 4
$ws.Range("F11").Value2 = 172
$ws.Range("H11").Value2 = "bedrooms"
$ws.Range("L11").Value2 = "stimuli/img_uxxo0.png"
$ws.Range("M11").Value2 = 71.74418604651163
$ws.Range("N11").Value2 = 48.44186046511628
$ws.Range("O11").Value2 = 60.09302325581395
$ws.Range("P11").Value2 = 43
$ws.Range("C12").Value2 = 4
$ws.Range("F12").Value2 = 173
$ws.Range("H12").Value2 = "kitchens"
$ws.Range("L12").Value2 = "stimuli/img_q577a.png"
$ws.Range("M12").Value2 = 81.26470588235294
$ws.Range("N12").Value2 = 59.08823529411764
$ws.Range("O12").Value2 = 70.17647058823529
$ws.Range("P12").Value2 = 34
$ws.Range("Q12").Value2 = 8
$ws.Range("R12").Value2 = 8
$ws.Range("S12").Value2 = 8
$ws.Range("C13").Value2 = 4
$ws.Range("F13").Value2 = 174
$ws.Range("H13").Value2 = "bedrooms"
$ws.Range("L13").Value2 = "stimuli/img_2js6m.png"
$ws.Range("M13").Value2 = 40.02777777777778
$ws.Range("N13").Value2 = 20.88888888888889
$ws.Range("O13").Value2 = 30.45833333333334
$ws.Range("P13").Value2 = 36
$ws.Range("C14").Value2 = 4
$ws.Range("F14").Value2 = 175
$ws.Range("H14").Value2 = "kitchens"
$ws.Range("L14").Value2 = "stimuli/img_a220l.png"
$ws.Range("M14").Value2 = 79.45945945945945
$ws.Range("N14").Value2 = 60.97297297297298
$ws.Range("O14").Value2 = 70.21621621621621
$ws.Range("P14").Value2 = 37
$ws.Range("Q14").Value2 = 8
$ws.Range("R14").Value2 = 8
$ws.Range("S14").Value2 = 8
$ws.Range("C15").Value2 = 4
$ws.Range("F15").Value2 = 176
$ws.Range("H15").Value2 = "kitchens"
$ws.Range("L15").Value2 = "stimuli/img_cv6mf.png"
$ws.Range("M15").Value2 = 66.8
$ws.Range("N15").Value2 = 42.08
$ws.Range("O15").Value2 = 54.44
$ws.Range("P15").Value2 = 25
$ws.Range("Q15").Value2 = 4
$ws.Range("R15").Value2 = 4
$ws.Range("S15").Value2 = 4
$ws.Range("C16").Value2 = 4
$ws.Range("F16").Value2 = 177
$ws.Range("H16").Value2 = "living_rooms"
$ws.Range("L16").Value2 = "stimuli/img_pbsj1.png"
$ws.Range("M16").Value2 = 73.88636363636364
$ws.Range("N16").Value2 = 51.52272727272727
$ws.Range("O16").Value2 = 62.70454545454545
$ws.Range("P16").Value2 = 44
$ws.Range("Q16").Value2 = 6
$ws.Range("R16").Value2 = 6
$ws.Range("S16").Value2 = 6
$ws.Range("C17").Value2 = 4
$ws.Range("F17").Value2 = 178
$ws.Range("H17").Value2 = "bedrooms"
$ws.Range("L17").Value2 = "stimuli/img_2pk6v.png"
$ws.Range("M17").Value2 = 85.08108108108108
$ws.Range("N17").Value2 = 66.16216216216216
$ws.Range("O17").Value2 = 75.62162162162161
$ws.Range("P17").Value2 = 37
$ws.Range("C18").Value2 = 4
$ws.Range("F18").Value2 = 179
$ws.Range("H18").Value2 = "bedrooms"
$ws.Range("L18").Value2 = "stimuli/img_a9acb.png"
$ws.Range("M18").Value2 = 77.11428571428571
$ws.Range("N18").Value2 = 58.42857142857143
$ws.Range("O18").Value2 = 67.77142857142857
$ws.Range("P18").Value2 = 35
$ws.Range("C19").Value2 = 4
$ws.Range("F19").Value2 = 180
$ws.Range("H19").Value2 = "bedrooms"
$ws.Range("L19").Value2 = "stimuli/img_ybbmx.png"
$ws.Range("M19").Value2 = 55.24324324324324
$ws.Range("N19").Value2 = 36.75675675675676
$ws.Range("O19").Value2 = 46
$ws.Range("P19").Value2 = 37
$ws.Range("C20").Value2 = 4
$ws.Range("F20").Value2 = 181
$ws.Range("H20").Value2 = "bedrooms"
$ws.Range("L20").Value2 = "stimuli/img_v8dra.png"
$ws.Range("M20").Value2 = 61.77272727272727
$ws.Range("N20").Value2 = 38.79545454545455
$ws.Range("O20").Value2 = 50.28409090909091
$ws.Range("P20").Value2 = 44
$ws.Range("C21").Value2 = 4
$ws.Range("F21").Value2 = 182
$ws.Range("H21").Value2 = "bedrooms"
$ws.Range("L21").Value2 = "stimuli/img_okvvw.png"
$ws.Range("M21").Value2 = 50.58333333333334
$ws.Range("N21").Value2 = 32.11111111111111
$ws.Range("O21").Value2 = 41.34722222222223
$ws.Range("P21").Value2 = 36
$ws.Range("C22").Value2 = 4
$ws.Range("F22").Value2 = 183
$ws.Range("H22").Value2 = "bedrooms"
$ws.Range("L22").Value2 = "stimuli/img_th7xh.png"
$ws.Range("M22").Value2 = 82.35897435897436
$ws.Range("N22").Value2 = 65.53846153846153
$ws.Range("O22").Value2 = 73.94871794871796
$ws.Range("P22").Value2 = 39
$ws.Range("C23").Value2 = 4
$ws.Range("F23").Value2 = 184
$ws.Range("H23").Value2 = "bedrooms"
$ws.Range("L23").Value2 = "stimuli/img_5m6x4.png"
$ws.Range("M23").Value2 = 80.23076923076923
$ws.Range("N23").Value2 = 58.41025641025641
$ws.Range("O23").Value2 = 69.32051282051282
$ws.Range("P23").Value2 = 39
$ws.Range("C24").Value2 = 4
$ws.Range("F24").Value2 = 185
$ws.Range("H24").Value2 = "bedrooms"
$ws.Range("L24").Value2 = "stimuli/img_71mhq.png"
$ws.Range("M24").Value2 = 69.34210526315789
$ws.Range("N24").Value2 = 47.02631578947368
$ws.Range("O24").Value2 = 58.18421052631579
$ws.Range("P24").Value2 = 38
$ws.Range("C25").Value2 = 4
$ws.Range("F25").Value2 = 186
$ws.Range("H25").Value2 = "bedrooms"
$ws.Range("L25").Value2 = "stimuli/img_fqgem.png"
$ws.Range("M25").Value2 = 80.75
$ws.Range("N25").Value2 = 61.475
$ws.Range("O25").Value2 = 71.1125
$ws.Range("P25").Value2 = 40
$ws.Range("C26").Value2 = 4
$ws.Range("F26").Value2 = 187
$ws.Range("H26").Value2 = "kitchens"
$ws.Range("L26").Value2 = "stimuli/img_2b8fp.png"
$ws.Range("M26").Value2 = 73.89189189189189
$ws.Range("N26").Value2 = 51.45945945945946
$ws.Range("O26").Value2 = 62.67567567567568
$ws.Range("P26").Value2 = 37
$ws.Range("Q26").Value2 = 6
$ws.Range("R26").Value2 = 6
$ws.Range("S26").Value2 = 6
$ws.Range("C27").Value2 = 4
$ws.Range("F27").Value2 = 188
$ws.Range("H27").Value2 = "bedrooms"
$ws.Range("L27").Value2 = "stimuli/img_wyctg.png"
$ws.Range("M27").Value2 = 33.44736842105263
$ws.Range("N27").Value2 = 11.39473684210526
$ws.Range("O27").Value2 = 22.42105263157895
$ws.Range("P27").Value2 = 38
